$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.676.50"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "2.656.26"
$ws.Range("E3").Value = "  -3.04%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "598.64"
$ws.Range("E5").Value = "  -1.60%  "
$ws.Range("D6").Value = "168.53"
$ws.Range("E6").Value = "  -1.21%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "2.656.32"
$ws.Range("E9").Value = "  -3.02%  "
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("E11").Value = "  +2.09%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "5.26"
$ws.Range("E13").Value = "  -2.03%  "
$ws.Range("D14").Value = "28.17"
$ws.Range("D15").Value = "3.137.88"
$ws.Range("E15").Value = "  -3.13%  "
$ws.Range("E16").Value = "  -3.52%  "
$ws.Range("D17").Value = "67.930.62"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "2.650.58"
$ws.Range("E18").Value = "  -2.93%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "8.35"
$ws.Range("E19").Value = "  +8.09%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "12.05"
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "364.21"
$ws.Range("D22").Value = "4.43"
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("E23").Value = "  -4.45%  "
$ws.Range("D24").Value = "11.03"
$ws.Range("E24").Value = "  +8.20%  "
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "70.97"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("D28").Value = "2.795.18"
$ws.Range("E28").Value = "  -2.79%  "
$ws.Range("E29").Value = "  -3.67%  "
$ws.Range("E30").Value = "  -0.25%  "
$ws.Range("D31").Value = "558.70"
$ws.Range("E31").Value = "  -6.18%  "
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("D34").Value = "1.94"
$ws.Range("E34").Value = "  -3.02%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -4.82%  "
$ws.Range("D38").Value = "158.28"
$ws.Range("E38").Value = "  -2.69%  "
$ws.Range("D39").Value = "19.46"
$ws.Range("E39").Value = "  -2.84%  "
$ws.Range("E40").Value = "  -2.67%  "
$ws.Range("E41").Value = "  -3.78%  "
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("D43").Value = "17.93"
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  -5.76%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "40.32"
$ws.Range("E46").Value = "  -2.04%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "0.600"
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0300"
$ws.Range("E48").Value = "  -4.20%  "
$ws.Range("D49").Value = "154.92"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("E51").Value = "  -3.18%  "
